$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '301.63'
    'E2' = '1.19%'
    'D3' = '31.81'
    'E3' = '1.68%'
    'D4' = '5.139'
    'E4' = '0.96%'
    'D5' = '0.07851'
    'E5' = '-2.45%'
    'D6' = '2.250'
    'E6' = '-11.01%'
    'D7' = '7.808'
    'E7' = '-0.12%'
    'D8' = '3.810'
    'E8' = '-0.15%'
    'D9' = '0.9283'
    'E9' = '1.14%'
    'D10' = '0.1775'
    'E10' = '2.65%'
    'D11' = '0.07655'
    'E11' = '4.40%'
    'D12' = '0.08861'
    'E12' = '2.41%'
    'D13' = '0.03103'
    'E13' = '2.49%'
    'E14' = '0.53%'
    'D15' = '0.001515'
    'E15' = '0.79%'
    'D16' = '0.005792'
    'E16' = '-3.78%'
    'D17' = '3.482'
    'E17' = '-0.41%'
    'D18' = '2.251'
    'E18' = '0.27%'
    'D20' = '0.1328'
    'E20' = '-0.70%'
    'D21' = '4.317'
    'E21' = '-5.89%'
    'E22' = '10.78%'
    'D23' = '0.04594'
    'E23' = '-0.32%'
    'E24' = '0.25%'
    'D25' = '0.004477'
    'E25' = '0.87%'
    'E26' = '4.04%'
    'E27' = '-1.30%'
    'D39' = '0.01781'
    'E39' = '-0.84%'
    'D40' = '0.04792'
    'E40' = '5.92%'
    'D41' = '0.007427'
    'E41' = '5.65%'
    'D42' = '0.1363'
    'E42' = '1.52%'
    'D43' = '0.002188'
    'E43' = '-2.35%'
    'D44' = '0.009819'
    'E44' = '-0.02%'
    'D45' = '0.00006248'
    'E45' = '-6.36%'
    'E46' = '-0.11%'
    'D48' = '0.7022'
    'E48' = '-14.42%'
    'E49' = '-0.11%'
    'E50' = '-0.11%'
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
    $range.Style = "Normal"
}

Write-Host "Updated $($updates.Count) cells"